# Auto-generated edit script applying the Goblin_Profits.xlsx diff
# Updates cached market-price / profit figures across 8 worksheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 83333510
$ws.Range("I6").Value = 125000100
$ws.Range("J6").Value = 339.5
$ws.Range("K6").Value = 375000300
$ws.Range("L6").Value = 1018.5
$ws.Range("M6").Value = -375000188
$ws.Range("N6").Value = -1242.5

# Row 26
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").ClearContents()

# Row 28
$ws.Range("H28").Value = 461.17648
$ws.Range("I28").Value = 466.07144
$ws.Range("K28").Value = 466.07144
$ws.Range("M28").Value = 18.92856

# Row 39
$ws.Range("H39").Value = 170.26666
$ws.Range("I39").Value = 183
$ws.Range("J39").Value = 119.333336
$ws.Range("K39").Value = 549
$ws.Range("L39").Value = 358.000008
$ws.Range("M39").Value = -253
$ws.Range("N39").Value = -950.000008

# Row 41
$ws.Range("H41").Value = 680.38464
$ws.Range("I41").Value = 640
$ws.Range("K41").Value = 640
$ws.Range("M41").Value = -200

# Row 74
$ws.Range("H74").Value = 3116.6
$ws.Range("I74").Value = 3017.7856
$ws.Range("K74").Value = 3017.7856
$ws.Range("M74").Value = -2081.7856

# Row 77
$ws.Range("H77").Value = 3116.6
$ws.Range("I77").Value = 3017.7856
$ws.Range("K77").Value = 15088.928
$ws.Range("M77").Value = -10408.928

# Row 80
$ws.Range("H80").Value = 2764.9565
$ws.Range("I80").Value = 851.7778
$ws.Range("J80").Value = 3994.8572
$ws.Range("K80").Value = 2555.3334
$ws.Range("L80").Value = 11984.5716
$ws.Range("M80").Value = -1557.3334
$ws.Range("N80").Value = -13980.5716

# Row 83
$ws.Range("H83").Value = 2764.9565
$ws.Range("I83").Value = 851.7778
$ws.Range("J83").Value = 3994.8572
$ws.Range("K83").Value = 7666.000199999999
$ws.Range("L83").Value = 35953.7148
$ws.Range("M83").Value = -2674.000199999999
$ws.Range("N83").Value = -45937.7148

# Row 136
$ws.Range("H136").Value = 136497.33
$ws.Range("J136").Value = 136497.33
$ws.Range("L136").Value = 136497.33
$ws.Range("N136").Value = -146697.33

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5591.7393
$ws.Range("I32").Value = 3702.8853
$ws.Range("J32").Value = 19994.25
$ws.Range("K32").Value = 3702.8853
$ws.Range("L32").Value = 19994.25
$ws.Range("M32").Value = -3415.8853
$ws.Range("N32").Value = -20568.25

# Row 102
$ws.Range("H102").Value = 6686.533
$ws.Range("I102").Value = 2833
$ws.Range("J102").Value = 9255.556
$ws.Range("K102").Value = 2833
$ws.Range("L102").Value = 9255.556
$ws.Range("M102").Value = -1211
$ws.Range("N102").Value = -12499.556

# Row 132
$ws.Range("H132").Value = 5527.6113
$ws.Range("I132").Value = 2961.3845
$ws.Range("K132").Value = 8884.1535
$ws.Range("M132").Value = -6354.1535

# Row 133
$ws.Range("H133").Value = 40984.25
$ws.Range("J133").Value = 40984.25
$ws.Range("L133").Value = 40984.25
$ws.Range("N133").Value = -46044.25

# Row 135
$ws.Range("H135").Value = 150000
$ws.Range("J135").Value = 150000
$ws.Range("L135").Value = 150000
$ws.Range("N135").Value = -160140

# Row 139
$ws.Range("H139").Value = 130715
$ws.Range("J139").Value = 130715
$ws.Range("L139").Value = 130715
$ws.Range("N139").Value = -140995

$ws = $wb.Worksheets.Item("CRP")
# Row 99
$ws.Range("H99").Value = 3734.5386
$ws.Range("I99").Value = 2112.9565
$ws.Range("K99").Value = 2112.9565
$ws.Range("M99").Value = -614.9564999999998

# Row 107
$ws.Range("H107").Value = 622.6087
$ws.Range("I107").Value = 647.65
$ws.Range("K107").Value = 647.65
$ws.Range("M107").Value = 1272.35

# Row 126
$ws.Range("H126").Value = 3734.5386
$ws.Range("I126").Value = 2112.9565
$ws.Range("K126").Value = 6338.869499999999
$ws.Range("M126").Value = -3868.869499999999

$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 579.6
$ws.Range("I5").Value = 589.7222
$ws.Range("J5").Value = 488.5
$ws.Range("K5").Value = 1769.1666
$ws.Range("L5").Value = 1465.5
$ws.Range("M5").Value = -1657.1666
$ws.Range("N5").Value = -1689.5

# Row 11
$ws.Range("H11").Value = 4791.7905
$ws.Range("I11").Value = 6603.2666
$ws.Range("J11").Value = 611.46155
$ws.Range("K11").Value = 19809.7998
$ws.Range("L11").Value = 1834.38465
$ws.Range("M11").Value = -19669.7998
$ws.Range("N11").Value = -2114.38465

# Row 135
$ws.Range("H135").Value = 579.6
$ws.Range("I135").Value = 589.7222
$ws.Range("J135").Value = 488.5
$ws.Range("K135").Value = 5307.499800000001
$ws.Range("L135").Value = 4396.5
$ws.Range("M135").Value = -2772.499800000001
$ws.Range("N135").Value = -9466.5

$ws = $wb.Worksheets.Item("GSM")
# Row 21
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

# Row 30
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()

# Row 64
$ws.Range("H64").Value = 500037500
$ws.Range("J64").Value = 500037500
$ws.Range("L64").Value = 500037500
$ws.Range("N64").Value = -500037996

# Row 67
$ws.Range("H67").Value = 500037500
$ws.Range("J67").Value = 500037500
$ws.Range("L67").Value = 500037500
$ws.Range("N67").Value = -500039216

# Row 107
$ws.Range("H107").Value = 1721.4762
$ws.Range("I107").Value = 699.9
$ws.Range("J107").Value = 2650.182
$ws.Range("K107").Value = 699.9
$ws.Range("L107").Value = 2650.182
$ws.Range("M107").Value = 1220.1
$ws.Range("N107").Value = -6490.182

# Row 132
$ws.Range("H132").Value = 7244.4546
$ws.Range("I132").Value = 2418.2
$ws.Range("K132").Value = 7254.599999999999
$ws.Range("M132").Value = -4724.599999999999

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 2265.889
$ws.Range("J40").Value = 2628.5715
$ws.Range("L40").Value = 2628.5715
$ws.Range("N40").Value = -2900.5715

# Row 46
$ws.Range("H46").Value = 2151.6
$ws.Range("I46").Value = 860
$ws.Range("J46").Value = 2474.5
$ws.Range("K46").Value = 860
$ws.Range("L46").Value = 2474.5
$ws.Range("M46").Value = -672
$ws.Range("N46").Value = -2850.5

# Row 55
$ws.Range("H55").Value = 1789.15
$ws.Range("I55").Value = 1292
$ws.Range("J55").Value = 2712.4285
$ws.Range("K55").Value = 1292
$ws.Range("L55").Value = 2712.4285
$ws.Range("M55").Value = -1119
$ws.Range("N55").Value = -3058.4285

# Row 125
$ws.Range("H125").Value = 30715
$ws.Range("J125").Value = 30715
$ws.Range("L125").Value = 30715
$ws.Range("N125").Value = -40555

$ws = $wb.Worksheets.Item("WVR")
# Row 46
$ws.Range("H46").Value = 57499.5
$ws.Range("J46").Value = 57499.5
$ws.Range("L46").Value = 57499.5
$ws.Range("N46").Value = -57961.5

# Row 64
$ws.Range("H64").Value = 185000
$ws.Range("J64").Value = 185000
$ws.Range("L64").Value = 185000
$ws.Range("N64").Value = -185496

# Row 67
$ws.Range("H67").Value = 185000
$ws.Range("J67").Value = 185000
$ws.Range("L67").Value = 185000
$ws.Range("N67").Value = -186716

# Row 74
$ws.Range("H74").Value = 24465.334
$ws.Range("J74").Value = 24832.875
$ws.Range("L74").Value = 24832.875
$ws.Range("N74").Value = -26704.875

# Row 77
$ws.Range("H77").Value = 24465.334
$ws.Range("J77").Value = 24832.875
$ws.Range("L77").Value = 74498.625
$ws.Range("N77").Value = -83858.625

# Row 81
$ws.Range("H81").Value = 5434.6665
$ws.Range("J81").Value = 6017
$ws.Range("L81").Value = 12034
$ws.Range("N81").Value = -14156

# Row 84
$ws.Range("H84").Value = 5434.6665
$ws.Range("J84").Value = 6017
$ws.Range("L84").Value = 60170
$ws.Range("N84").Value = -70778

# Row 107
$ws.Range("H107").Value = 519.2963
$ws.Range("J107").Value = 588.44446
$ws.Range("L107").Value = 1765.33338
$ws.Range("N107").Value = -5605.33338

# Row 133
$ws.Range("H133").Value = 52497
$ws.Range("J133").Value = 52497
$ws.Range("L133").Value = 52497
$ws.Range("N133").Value = -62617

# Row 134
$ws.Range("H134").Value = 57499.5
$ws.Range("J134").Value = 57499.5
$ws.Range("L134").Value = 172498.5
$ws.Range("N134").Value = -177568.5
